# valid-utilisation-report.xlsx — fn 997: file type validation & virus scanning
# Fix the "Facility utilisatin" header typo, tweak one figure, move the
# selection, and shrink the header row height to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the misspelled "Facility utilisatin" column header (F1) -> "Facility utilisation".
$ws.Range("F1").Value = "Facility utilisation"

# Row 4's utilisation percentage was corrected from 124.758 to 124.75.
$ws.Range("G4").Value = 124.75

# Shrink header row 1's height now that the corrected heading wraps differently.
$ws.Rows.Item(1).RowHeight = 87.5

# Move the active selection to G14.
[void]$ws.Range("G14").Select()
